$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns are treated as text so values like "1.016"
# or percentages are stored literally, matching the source inline strings
# rather than being auto-converted into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '27.967.36'
$ws.Cells.Item(2, 5).Value = '  +0.10%  '

$ws.Cells.Item(3, 4).Value = '1.845.81'
$ws.Cells.Item(3, 5).Value = '  -0.77%  '

$ws.Cells.Item(4, 4).Value = '1.016'
$ws.Cells.Item(4, 5).Value = '  +1.16%  '

$ws.Cells.Item(5, 4).Value = '331.60'
$ws.Cells.Item(5, 5).Value = '  -1.28%  '

$ws.Cells.Item(6, 4).Value = '1.016'
$ws.Cells.Item(6, 5).Value = '  +1.34%  '

$ws.Cells.Item(7, 4).Value = '0.4498'
$ws.Cells.Item(7, 5).Value = '  -4.07%  '

$ws.Cells.Item(8, 4).Value = '0.3884'
$ws.Cells.Item(8, 5).Value = '  +0.26%  '

$ws.Cells.Item(9, 4).Value = '47.81'
$ws.Cells.Item(9, 5).Value = '  +2.12%  '

$ws.Cells.Item(10, 4).Value = '0.07738'
$ws.Cells.Item(10, 5).Value = '  -2.75%  '

$ws.Cells.Item(11, 4).Value = '0.9730'
$ws.Cells.Item(11, 5).Value = '  -0.17%  '

$ws.Cells.Item(12, 4).Value = '21.10'
$ws.Cells.Item(12, 5).Value = '  -1.43%  '

$ws.Cells.Item(13, 4).Value = '1.853.64'
$ws.Cells.Item(13, 5).Value = '  +0.53%  '

$ws.Cells.Item(14, 4).Value = '5.774'
$ws.Cells.Item(14, 5).Value = '  -2.31%  '

$ws.Cells.Item(15, 4).Value = '6.935'
$ws.Cells.Item(15, 5).Value = '  -3.47%  '

$ws.Cells.Item(16, 4).Value = '1.020'
$ws.Cells.Item(16, 5).Value = '  +1.61%  '

$ws.Cells.Item(17, 4).Value = '87.05'
$ws.Cells.Item(17, 5).Value = '  -4.57%  '

$ws.Cells.Item(18, 4).Value = '0.06569'
$ws.Cells.Item(18, 5).Value = '  -0.84%  '

$ws.Cells.Item(19, 4).Value = '0.00001013'
$ws.Cells.Item(19, 5).Value = '  -1.82%  '

$ws.Cells.Item(20, 4).Value = '16.85'
$ws.Cells.Item(20, 5).Value = '  -2.98%  '

$ws.Cells.Item(21, 4).Value = '1.006'
$ws.Cells.Item(21, 5).Value = '  +0.11%  '

$ws.Cells.Item(22, 4).Value = '27.932.45'
$ws.Cells.Item(22, 5).Value = '  +0.02%  '

$ws.Cells.Item(23, 4).Value = '5.294'
$ws.Cells.Item(23, 5).Value = '  -1.44%  '

$ws.Cells.Item(24, 4).Value = '10.57'
$ws.Cells.Item(24, 5).Value = '  -2.65%  '

$ws.Cells.Item(25, 4).Value = '2.262'
$ws.Cells.Item(25, 5).Value = '  -1.48%  '

$ws.Cells.Item(26, 4).Value = '2.088.25'
$ws.Cells.Item(26, 5).Value = '  +0.73%  '

$ws.Cells.Item(27, 4).Value = '155.22'
$ws.Cells.Item(27, 5).Value = '  -2.40%  '

$ws.Cells.Item(28, 4).Value = '19.13'
$ws.Cells.Item(28, 5).Value = '  -2.00%  '

$ws.Cells.Item(29, 4).Value = '2.026'
$ws.Cells.Item(29, 5).Value = '  -2.97%  '

$ws.Cells.Item(30, 4).Value = '5.228'
$ws.Cells.Item(30, 5).Value = '  -3.56%  '

$ws.Cells.Item(31, 4).Value = '115.99'
$ws.Cells.Item(31, 5).Value = '  -2.63%  '

$ws.Cells.Item(32, 2).Value = 'Stellar'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(32, 4).Value = '0.09273'
$ws.Cells.Item(32, 5).Value = '  -2.04%  '

$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(33, 4).Value = '0.9297'
$ws.Cells.Item(33, 5).Value = '  -2.13%  '

$ws.Cells.Item(34, 4).Value = '3.636'
$ws.Cells.Item(34, 5).Value = '  +1.37%  '

$ws.Cells.Item(35, 4).Value = '1.366'
$ws.Cells.Item(35, 5).Value = '  +1.89%  '

$ws.Cells.Item(36, 4).Value = '5.181'
$ws.Cells.Item(36, 5).Value = '  -1.91%  '

$ws.Cells.Item(37, 4).Value = '0.05987'
$ws.Cells.Item(37, 5).Value = '  -1.42%  '

$ws.Cells.Item(38, 4).Value = '0.02174'
$ws.Cells.Item(38, 5).Value = '  -2.20%  '

$ws.Cells.Item(39, 4).Value = '8.144'
$ws.Cells.Item(39, 5).Value = '  -1.35%  '

$ws.Cells.Item(40, 2).Value = 'Frax'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(40, 4).Value = '1.016'
$ws.Cells.Item(40, 5).Value = '  +1.31%  '

$ws.Cells.Item(41, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(41, 4).Value = '1.143'
$ws.Cells.Item(41, 5).Value = '  -1.57%  '

$ws.Cells.Item(42, 4).Value = '0.5645'
$ws.Cells.Item(42, 5).Value = '  -3.51%  '

$ws.Cells.Item(43, 2).Value = 'Algorand'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(43, 4).Value = '0.1783'
$ws.Cells.Item(43, 5).Value = '  -4.04%  '

$ws.Cells.Item(44, 2).Value = 'Aptos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(44, 4).Value = '9.847'
$ws.Cells.Item(44, 5).Value = '  -2.89%  '

$ws.Cells.Item(45, 4).Value = '1.230'
$ws.Cells.Item(45, 5).Value = '  -4.15%  '

$ws.Cells.Item(46, 4).Value = '2.239'
$ws.Cells.Item(46, 5).Value = '  +23.10%  '

$ws.Cells.Item(47, 2).Value = 'Cronos'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(47, 4).Value = '0.07170'
$ws.Cells.Item(47, 5).Value = '  +4.40%  '

$ws.Cells.Item(48, 2).Value = 'Decentraland'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(48, 4).Value = '0.5348'
$ws.Cells.Item(48, 5).Value = '  -2.52%  '

$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).Value = '11.59'
$ws.Cells.Item(49, 5).Value = '  -3.72%  '

$ws.Cells.Item(50, 4).Value = '1.856'
$ws.Cells.Item(50, 5).Value = '  -4.49%  '

$ws.Cells.Item(51, 4).Value = '109.18'
$ws.Cells.Item(51, 5).Value = '  -2.13%  '
